# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# 1) Updates the "Valor Mora" total (E11) and "Cant. Periodos" (F13).
# 2) Corrects the second worker's "Salario Basico" (G17) to match the
#    other rows.
# 3) Adds a second "Periodo Mora" (2508) block for the first two workers,
#    inserting two new table rows (19 and 20) right after the existing
#    data rows, re-using the same look (fonts/fills/borders) as the rest
#    of the table: the regular middle-row style for the non-last rows and
#    the thicker bottom-border style for the new last row.
# 4) The signature block (rows 23-24) is naturally pushed down to 25-26 by
#    the row insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary values ---
$ws.Range("E11").Value = 284700
$ws.Range("F13").Value = 2

# --- Fix existing row 17 "Salario Basico" value ---
$ws.Range("G17").Value = 1423500

# --- Insert two new rows right after row 18 (they become rows 19 and 20) ---
$ws.Rows.Item(19).EntireRow.Insert()
$ws.Rows.Item(19).EntireRow.Insert()

# Row 18 (source of the "last row" look, thicker bottom border) is copied
# onto the new final row (20) BEFORE row 18's own formatting is changed.
$ws.Range("B18:J18").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# Row 18 becomes a regular (non-last) row, matching row 17's look.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# Row 19 is also a regular (non-last) row.
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# --- Row 19 content: worker 1 (RAMON HUMBERTO DOMINGUEZ CANTILLO), period 2508 ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73117838"
$ws.Range("D19").Value = "RAMON HUMBERTO DOMINGUEZ CANTILLO"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# --- Row 20 content: worker 2 (GLAYDIS MILENA GALVIS GARCIA), period 2508 ---
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1104377019"
$ws.Range("D20").Value = "GLAYDIS MILENA GALVIS GARCIA"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500
